$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 2100.5715
$ws.Range("J17").Value = 2100.5715
$ws.Range("L17").Value = 6301.7145
$ws.Range("N17").Value = -6637.7145
$ws.Range("H28").Value = 60689.47
$ws.Range("I28").Value = 68712.2
$ws.Range("K28").Value = 68712.2
$ws.Range("M28").Value = -68227.2
$ws.Range("H76").Value = 71435660
$ws.Range("J76").Value = 142863710
$ws.Range("L76").Value = 142863710
$ws.Range("N76").Value = -142864340
$ws.Range("H79").Value = 71435660
$ws.Range("J79").Value = 142863710
$ws.Range("L79").Value = 142863710
$ws.Range("N79").Value = -142865894
$ws.Range("H86").Value = 1951845.1
$ws.Range("I86").Value = 1371.375
$ws.Range("J86").Value = 4788897.5
$ws.Range("K86").Value = 1371.375
$ws.Range("L86").Value = 4788897.5
$ws.Range("M86").Value = -248.375
$ws.Range("N86").Value = -4791143.5
$ws.Range("H89").Value = 1951845.1
$ws.Range("I89").Value = 1371.375
$ws.Range("J89").Value = 4788897.5
$ws.Range("K89").Value = 6856.875
$ws.Range("L89").Value = 23944487.5
$ws.Range("M89").Value = -1240.875
$ws.Range("N89").Value = -23955719.5
$ws.Range("H98").Value = 747.2258
$ws.Range("I98").Value = 623.5862
$ws.Range("K98").Value = 623.5862
$ws.Range("M98").Value = 874.4138
$ws.Range("H122").Value = 747.2258
$ws.Range("I122").Value = 623.5862
$ws.Range("K122").Value = 1870.7586
$ws.Range("M122").Value = 579.2414000000001
$ws.Range("H132").Value = 2207.878
$ws.Range("I132").Value = 2236.7778
$ws.Range("K132").Value = 6710.3334
$ws.Range("M132").Value = -4180.3334
$ws.Range("H135").Value = 2021.4839
$ws.Range("I135").Value = 1842
$ws.Range("J135").Value = 3233
$ws.Range("K135").Value = 16578
$ws.Range("L135").Value = 29097
$ws.Range("M135").Value = -14043
$ws.Range("N135").Value = -34167
$ws.Range("H137").Value = 1677.541
$ws.Range("I137").Value = 1380.1951
$ws.Range("K137").Value = 4140.5853
$ws.Range("M137").Value = -1590.5853

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 5468.375
$ws.Range("I32").Value = 5004.067
$ws.Range("J32").Value = 12433
$ws.Range("K32").Value = 5004.067
$ws.Range("L32").Value = 12433
$ws.Range("M32").Value = -4717.067
$ws.Range("N32").Value = -13007
$ws.Range("H74").Value = 2085.5715
$ws.Range("I74").Value = 2268.5908
$ws.Range("J74").Value = 1775.8462
$ws.Range("K74").Value = 2268.5908
$ws.Range("L74").Value = 1775.8462
$ws.Range("M74").Value = -1394.5908
$ws.Range("N74").Value = -3523.8462
$ws.Range("H77").Value = 2085.5715
$ws.Range("I77").Value = 2268.5908
$ws.Range("J77").Value = 1775.8462
$ws.Range("K77").Value = 11342.954
$ws.Range("L77").Value = 8879.231
$ws.Range("M77").Value = -6974.954
$ws.Range("N77").Value = -17615.231
$ws.Range("H110").Value = 202109.08
$ws.Range("I110").Value = 265429.94
$ws.Range("K110").Value = 265429.94
$ws.Range("M110").Value = -263384.94
$ws.Range("H132").Value = 6536
$ws.Range("I132").Value = 5451.7
$ws.Range("J132").Value = 9246.75
$ws.Range("K132").Value = 16355.1
$ws.Range("L132").Value = 27740.25
$ws.Range("M132").Value = -13825.1
$ws.Range("N132").Value = -32800.25

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H75").Value = 19249.75
$ws.Range("I75").Value = 19249.75
$ws.Range("K75").Value = 19249.75
$ws.Range("M75").Value = -18313.75
$ws.Range("H78").Value = 19249.75
$ws.Range("I78").Value = 19249.75
$ws.Range("K78").Value = 57749.25
$ws.Range("M78").Value = -53069.25
$ws.Range("H82").Value = 34771.4
$ws.Range("I82").Value = 5416.2856
$ws.Range("J82").Value = 103266.664
$ws.Range("K82").Value = 5416.2856
$ws.Range("L82").Value = 103266.664
$ws.Range("M82").Value = -5033.2856
$ws.Range("N82").Value = -104032.664
$ws.Range("H85").Value = 34771.4
$ws.Range("I85").Value = 5416.2856
$ws.Range("J85").Value = 103266.664
$ws.Range("K85").Value = 5416.2856
$ws.Range("L85").Value = 103266.664
$ws.Range("M85").Value = -4090.2856
$ws.Range("N85").Value = -105918.664
$ws.Range("H134").Value = 25936.223
$ws.Range("I134").Value = 3674.3635
$ws.Range("J134").Value = 87156.336
$ws.Range("K134").Value = 11023.0905
$ws.Range("L134").Value = 261469.008
$ws.Range("M134").Value = -8488.0905
$ws.Range("N134").Value = -266539.008

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 4394.9414
$ws.Range("I16").Value = 4275
$ws.Range("J16").Value = 4682.8
$ws.Range("K16").Value = 4275
$ws.Range("L16").Value = 4682.8
$ws.Range("M16").Value = -3988
$ws.Range("N16").Value = -5256.8
$ws.Range("H105").Value = 1152.125
$ws.Range("I105").Value = 1084.3334
$ws.Range("J105").Value = 1355.5
$ws.Range("K105").Value = 1084.3334
$ws.Range("L105").Value = 1355.5
$ws.Range("M105").Value = 662.6666
$ws.Range("N105").Value = -4849.5
$ws.Range("H113").Value = 4394.9414
$ws.Range("I113").Value = 4275
$ws.Range("J113").Value = 4682.8
$ws.Range("K113").Value = 4275
$ws.Range("L113").Value = 4682.8
$ws.Range("M113").Value = -2105
$ws.Range("N113").Value = -9022.799999999999
$ws.Range("H122").Value = 2089.7368
$ws.Range("I122").Value = 2120.5334
$ws.Range("J122").Value = 1974.25
$ws.Range("K122").Value = 6361.600199999999
$ws.Range("L122").Value = 5922.75
$ws.Range("M122").Value = -3911.600199999999
$ws.Range("N122").Value = -10822.75
$ws.Range("H134").Value = 288536.84
$ws.Range("I134").Value = 2960.394
$ws.Range("J134").Value = 5000548.5
$ws.Range("K134").Value = 8881.181999999999
$ws.Range("L134").Value = 15001645.5
$ws.Range("M134").Value = -6346.181999999999
$ws.Range("N134").Value = -15006715.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 450837.34
$ws.Range("I5").Value = 38656.617
$ws.Range("K5").Value = 115969.851
$ws.Range("M5").Value = -115857.851
$ws.Range("H86").Value = 500.46155
$ws.Range("I86").Value = 235.77777
$ws.Range("J86").Value = 1096
$ws.Range("K86").Value = 707.33331
$ws.Range("L86").Value = 3288
$ws.Range("M86").Value = 478.66669
$ws.Range("N86").Value = -5660
$ws.Range("H89").Value = 500.46155
$ws.Range("I89").Value = 235.77777
$ws.Range("J89").Value = 1096
$ws.Range("K89").Value = 2121.99993
$ws.Range("L89").Value = 9864
$ws.Range("M89").Value = 3806.00007
$ws.Range("N89").Value = -21720
$ws.Range("H97").Value = 553.5714
$ws.Range("I97").Value = 587.5
$ws.Range("J97").Value = 350
$ws.Range("K97").Value = 1762.5
$ws.Range("L97").Value = 1050
$ws.Range("M97").Value = -1266.5
$ws.Range("N97").Value = -2042
$ws.Range("H113").Value = 1612232.8
$ws.Range("I113").Value = 4116350
$ws.Range("J113").Value = 2443.0715
$ws.Range("K113").Value = 12349050
$ws.Range("L113").Value = 7329.2145
$ws.Range("M113").Value = -12346880
$ws.Range("N113").Value = -11669.2145
$ws.Range("H122").Value = 2261.9524
$ws.Range("J122").Value = 2921.5833
$ws.Range("L122").Value = 26294.2497
$ws.Range("N122").Value = -31194.2497
$ws.Range("H132").Value = 503328.3
$ws.Range("J132").Value = 670694.9399999999
$ws.Range("L132").Value = 6036254.459999999
$ws.Range("N132").Value = -6041314.459999999
$ws.Range("H135").Value = 450837.34
$ws.Range("I135").Value = 38656.617
$ws.Range("K135").Value = 347909.553
$ws.Range("M135").Value = -345374.553

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 4118
$ws.Range("J15").Value = 4118
$ws.Range("L15").Value = 4118
$ws.Range("N15").Value = -4694
$ws.Range("H81").Value = 4118
$ws.Range("J81").Value = 4118
$ws.Range("L81").Value = 4118
$ws.Range("N81").Value = -6114
$ws.Range("H84").Value = 4118
$ws.Range("J84").Value = 4118
$ws.Range("L84").Value = 12354
$ws.Range("N84").Value = -22338
$ws.Range("H102").Value = 2753.1333
$ws.Range("I102").Value = 1456.1052
$ws.Range("K102").Value = 1456.1052
$ws.Range("M102").Value = 165.8948
$ws.Range("H132").Value = 43075.184
$ws.Range("J132").Value = 168633.33
$ws.Range("L132").Value = 505899.99
$ws.Range("N132").Value = -510959.99
$ws.Range("H139").Value = 0
$ws.Range("J139").Value = 0
$ws.Range("L139").Value = 0
$ws.Range("N139").ClearContents()

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H40").Value = 7999.778
$ws.Range("I40").Value = 7998
$ws.Range("J40").Value = 8000
$ws.Range("K40").Value = 7998
$ws.Range("L40").Value = 8000
$ws.Range("M40").Value = -7862
$ws.Range("N40").Value = -8272
$ws.Range("H75").Value = 34500
$ws.Range("J75").Value = 34500
$ws.Range("L75").Value = 34500
$ws.Range("N75").Value = -36372
$ws.Range("H78").Value = 34500
$ws.Range("J78").Value = 34500
$ws.Range("L78").Value = 103500
$ws.Range("N78").Value = -112860
$ws.Range("H122").Value = 4993.884
$ws.Range("I122").Value = 4387.1113
$ws.Range("K122").Value = 13161.3339
$ws.Range("M122").Value = -10711.3339

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 27029632
$ws.Range("I122").Value = 37039108
$ws.Range("J122").Value = 4043.3
$ws.Range("K122").Value = 111117324
$ws.Range("L122").Value = 12129.9
$ws.Range("M122").Value = -111114874
$ws.Range("N122").Value = -17029.9
$ws.Range("H123").Value = 75000
$ws.Range("J123").Value = 75000
$ws.Range("L123").Value = 75000
$ws.Range("N123").Value = -84800
$ws.Range("H132").Value = 51226.227
$ws.Range("I132").Value = 5448
$ws.Range("K132").Value = 16344
$ws.Range("M132").Value = -13814
$ws.Range("H136").Value = 51642.348
$ws.Range("I136").Value = 13299.132
$ws.Range("J136").Value = 184100.73
$ws.Range("K136").Value = 39897.396
$ws.Range("L136").Value = 552302.1900000001
$ws.Range("M136").Value = -37347.396
$ws.Range("N136").Value = -557402.1900000001
